$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7:8 hold the "a1"/"b2" sample pair whose formatting (bold, bordered,
# centered style applied to column A) should be reused for the two new
# rows being appended at the bottom of the log (rows 11:12).
$ws.Range("A7:K8").Copy()
$ws.Range("A11:K12").PasteSpecial(-4122)   # xlPasteFormats

# Row 11 ("a1")
$ws.Range("A11").Value = "a1"
$ws.Range("B11").Value = 0.8752631545066833
$ws.Range("C11").Value = 0.5897498726844788
$ws.Range("D11").Value = 0.837644100189209
$ws.Range("E11").Value = 0.8836870193481445
$ws.Range("F11").Value = 0.8707107901573181
$ws.Range("G11").Value = 148.128173828125
$ws.Range("H11").Value = 13.14312839508057
$ws.Range("I11").Value = 13.35710144042969
$ws.Range("J11").Value = 86.67254638671875
$ws.Range("K11").Value = 91.16632843017578

# Row 12 ("b2")
$ws.Range("A12").Value = "b2"
$ws.Range("B12").Value = 0.875263512134552
$ws.Range("C12").Value = 0.5897493362426758
$ws.Range("D12").Value = 0.8376448750495911
$ws.Range("E12").Value = 0.8836870193481445
$ws.Range("F12").Value = 0.8707107305526733
$ws.Range("G12").Value = 148.1097259521484
$ws.Range("H12").Value = 13.14150238037109
$ws.Range("I12").Value = 13.35542774200439
$ws.Range("J12").Value = 86.67243957519531
$ws.Range("K12").Value = 91.16508483886719
